# Edit the "Seat Assignments" sheet:
#  - Swap the Original Position (A) and ContestantID (C) values between
#    row 3 and row 4 (the two contestants effectively swapped seats).
#  - Update the Seat column (E) values for rows 2-4 to reflect the
#    assigned seat for each row (B2, B3, B4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seat Assignments")

# Capture existing values for row 3 and row 4 before overwriting them
# (use Value2 for reading - Value getter in this runtime does not
# reliably return the underlying cell value)
$a3 = $ws.Range("A3").Value2
$c3 = $ws.Range("C3").Value2
$a4 = $ws.Range("A4").Value2
$c4 = $ws.Range("C4").Value2

# Swap "Original Position" (A) values between rows 3 and 4
$ws.Range("A3").Value = $a4
$ws.Range("A4").Value = $a3

# Swap "ContestantID" (C) values between rows 3 and 4
$ws.Range("C3").Value = $c4
$ws.Range("C4").Value = $c3

# Update the Seat column (E) for rows 2-4
$ws.Range("E2").Value = "B2"
$ws.Range("E3").Value = "B3"
$ws.Range("E4").Value = "B4"
